$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '43.932.39'
$c.Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  +0.37%  '
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '2.358.25'
$c.Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +0.45%  '
$ws.Cells.Item(4, 5).Value = '  +0.17%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '0.671'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +3.34%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '235.73'
$c.Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +1.28%  '
$ws.Cells.Item(7, 5).Value = '  +11.37%  '
$ws.Cells.Item(8, 5).Value = '  -0.06%  '
$ws.Cells.Item(9, 5).Value = '  +24.74%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '0.0987'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +1.47%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '28.08'
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +4.39%  '
$ws.Cells.Item(12, 5).Value = '  +2.14%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '2.706.12'
$c.Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +0.36%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '16.82'
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +8.72%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '6.80'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +10.30%  '
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '0.890'
$c.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +4.52%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '2.375.92'
$c.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +1.37%  '
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '43.901.06'
$c.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +0.47%  '
$ws.Cells.Item(19, 5).Value = '  +2.50%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '78.39'
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +5.83%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '6.41'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +2.16%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '253.56'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +1.49%  '
$ws.Cells.Item(23, 5).Value = '  -0.06%  '
$ws.Cells.Item(24, 5).Value = '  -1.09%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '2.50'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +2.98%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '10.68'
$c.Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +7.43%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '2.29'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.23%  '
$ws.Cells.Item(28, 5).Value = '  +0.63%  '
$ws.Cells.Item(29, 2).Value = 'ImmutableX'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '1.59'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +9.48%  '
$ws.Cells.Item(30, 2).Value = 'Monero'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = '@'
$c.Value = '172.59'
$c.Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -1.36%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '0.130'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -0.36%  '
$ws.Cells.Item(32, 5).Value = '  +5.21%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '5.20'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +4.19%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = '@'
$c.Value = '0.0721'
$c.Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +4.60%  '
$ws.Cells.Item(35, 5).Value = '  +4.85%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '3.77'
$c.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +1.01%  '
$ws.Cells.Item(37, 5).Value = '  -0.66%  '
$ws.Cells.Item(38, 5).Value = '  -2.22%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '0.0272'
$c.Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +7.25%  '
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = '@'
$c.Value = '19.23'
$c.Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +7.57%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '8.98'
$c.Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -2.54%  '
$ws.Cells.Item(42, 5).Value = '  -0.10%  '
$ws.Cells.Item(43, 2).Value = 'ARBITRUM'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '1.17'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -1.48%  '
$ws.Cells.Item(44, 2).Value = 'Cronos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '0.0979'
$c.Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +2.30%  '
$ws.Cells.Item(45, 2).Value = 'Algorand'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '0.186'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +14.94%  '
$ws.Cells.Item(46, 5).Value = '  +2.04%  '
$ws.Cells.Item(47, 2).Value = 'FTXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '4.45'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +2.65%  '
$ws.Cells.Item(48, 2).Value = 'Aave'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '98.20'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -1.36%  '
$ws.Cells.Item(49, 2).Value = 'Maker'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '1.435.88'
$c.Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -0.89%  '
$ws.Cells.Item(50, 2).Value = 'NEARProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '2.31'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +0.34%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '2.582.06'
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +0.43%  '
